$d = $word.ActiveDocument

$pairs = @(
    @("819÷7=", "661÷2="),
    @("445÷9=", "250÷6="),
    @("576÷4=", "408÷6="),
    @("785÷2=", "476÷2="),
    @("975÷5=", "155÷7="),
    @("434÷3=", "786÷9="),
    @("855÷7=", "300÷4="),
    @("475÷9=", "139÷4="),
    @("872÷7=", "467÷5="),
    @("946÷9=", "285÷4="),
    @("498÷6=", "306÷9="),
    @("500÷9=", "188÷6="),
    @("185÷3=", "468÷9="),
    @("186÷2=", "172÷7="),
    @("672÷9=", "912÷7="),
    @("437÷7=", "374÷5="),
    @("290÷5=", "821÷6="),
    @("471÷7=", "663÷5="),
    @("633÷6=", "238÷9="),
    @("154÷6=", "691÷8="),
    @("750÷2=", "279÷7="),
    @("824÷9=", "478÷4="),
    @("844÷2=", "373÷7="),
    @("360÷7=", "434÷7="),
    @("321÷3=", "583÷4=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
